# Regenerate orders with updated distance/size codes.
# The experiment's distance conditions and one size condition were
# renamed (D64->D69, D80->D86, D51->D55, S30->S31). These codes appear
# embedded inside many strings throughout the sheet (Condition names,
# left/right stimulus filenames, and the standalone Distance/Size lookup
# values), so apply the substitutions as whole-token replacements across
# every used cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange

# Order matters only in that each pass is independent/non-overlapping,
# so any order works; keep it in the same order the renames were
# introduced for readability.
$used.Replace("D64", "D69") | Out-Null
$used.Replace("D80", "D86") | Out-Null
$used.Replace("D51", "D55") | Out-Null
$used.Replace("S30", "S31") | Out-Null
